# Apply the changes described by the commit diff:
#  1. Update the "Alphabets" config value in cell B6 from "ab,ac" to "ab,ca"
#  2. Update the active cell selection on the sheet from B18 to B7
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "ab,ca"

$ws.Range("B7").Select()
